# Mise à jour de certains champs de Modules et de Professeurs
# - Ajout d'une colonne "Matières enseignés" (E1)
# - Largeurs de colonnes C, D, E
# - Sélection active sur E6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Matières enseignés" column
$ws.Range("E1").Value = "Matières enseignés"

# Column widths (character units) for C, D, E
$ws.Columns("C").ColumnWidth = 27.5703125
$ws.Columns("D").ColumnWidth = 15.7109375
$ws.Columns("E").ColumnWidth = 31.7109375

# Move/restore the active selection
[void]$ws.Range("E6").Select()
